{"js": "// Fill in the empty \"Hex\" and \"Binary\" columns of the random-number\n// table (rows for numbers 1-15) with their computed values.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// rowIndex -> [hex, binary] ; rows are 0-indexed and row 0 is the header,\n// so table row 1 is the data row for \"Random Number\" == 0 (already filled),\n// table row 2 is the data row for \"1\", etc.\nconst rows = [\n  [\"B244\", \"1011 0010 0100 0100\"],\n  [\"0047\", \"0000 0000 0100 0111\"],\n  [\"0C54\", \"0000 11001010 0100\"],\n  [\"38C6\", \"0011 1000 1100 0110\"],\n  [\"BEF5\", \"1011 1110 1111 0101\"],\n  [\"BDCA\", \"1011 1101 1100 1010\"],\n  [\"251F\", \"0010 0101 0001 1111\"],\n  [\"316B\", \"0011 0001 0110 1011\"],\n  [\"C138\", \"1100 0001 0011 1000\"],\n  [\"5368\", \"0101 0011 0110 1000\"],\n  [\"8B71\", \"1000 1011 0111 0001\"],\n  [\"8DA9\", \"1000 1101 1010 1001\"],\n  [\"E5CA\", \"1110 0101 1010 1100\"],\n  [\"F5F9\", \"1111 0101 1111 1001 \"],\n];\n\nlet rowIndex = 2; // first empty data row (Random Number == 1)\nfor (const [hex, binary] of rows) {\n  const hexCell = table.getCell(rowIndex, 1);\n  hexCell.body.insertText(hex, Word.InsertLocation.replace);\n\n  const binaryCell = table.getCell(rowIndex, 2);\n  binaryCell.body.insertText(binary, Word.InsertLocation.replace);\n\n  rowIndex++;\n}\nawait context.sync();\n\n// Last row (Random Number == 15) has its binary value split across four\n// separate runs, joined by three single-space runs, instead of one run.\nconst lastRow = rowIndex; // 16\nconst lastHexCell = table.getCell(lastRow, 1);\nlastHexCell.body.insertText(\"7A67\", Word.InsertLocation.replace);\n\nconst lastBinaryCell = table.getCell(lastRow, 2);\nconst ooxml = `<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>0111</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:r><w:t>1010</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:r><w:t>0110</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:r><w:t>0111</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\nlastBinaryCell.body.insertOoxml(ooxml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Fill in the empty \"Hex\" and \"Binary\" columns of the random-number\n# table (rows for numbers 1-15) with their computed values.\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nfunction Set-CellRuns($cell, [string[]]$texts) {\n    $runsXml = \"\"\n    foreach ($t in $texts) {\n        $escaped = $t -replace '&','&amp;' -replace '<','&lt;' -replace '>','&gt;'\n        if ($t -match '^\\s|\\s$|^$') {\n            $runsXml += \"<w:r><w:t xml:space=`\"preserve`\">$escaped</w:t></w:r>\"\n        } else {\n            $runsXml += \"<w:r><w:t>$escaped</w:t></w:r>\"\n        }\n    }\n    $ooxml = '<?xml version=\"1.0\" encoding=\"utf-8\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    [void]$cell.Range.InsertXML($ooxml)\n}\n\n# table row index (1-based) -> [hex, binary] ; row 1 is the header and\n# row 2 is \"Random Number\" == 0 (already filled), so row 3 is \"1\", ... row 17 is \"15\".\n$rows = @(\n    ,@(\"B244\", \"1011 0010 0100 0100\")\n    ,@(\"0047\", \"0000 0000 0100 0111\")\n    ,@(\"0C54\", \"0000 11001010 0100\")\n    ,@(\"38C6\", \"0011 1000 1100 0110\")\n    ,@(\"BEF5\", \"1011 1110 1111 0101\")\n    ,@(\"BDCA\", \"1011 1101 1100 1010\")\n    ,@(\"251F\", \"0010 0101 0001 1111\")\n    ,@(\"316B\", \"0011 0001 0110 1011\")\n    ,@(\"C138\", \"1100 0001 0011 1000\")\n    ,@(\"5368\", \"0101 0011 0110 1000\")\n    ,@(\"8B71\", \"1000 1011 0111 0001\")\n    ,@(\"8DA9\", \"1000 1101 1010 1001\")\n    ,@(\"E5CA\", \"1110 0101 1010 1100\")\n    ,@(\"F5F9\", \"1111 0101 1111 1001 \")\n)\n\n$rowIndex = 3\nforeach ($pair in $rows) {\n    $hexCell = $table.Cell($rowIndex, 2)\n    Set-CellRuns $hexCell @($pair[0])\n\n    $binaryCell = $table.Cell($rowIndex, 3)\n    Set-CellRuns $binaryCell @($pair[1])\n\n    $rowIndex++\n}\n\n# Last row (Random Number == 15) has its binary value split across four\n# separate runs, joined by three single-space runs, instead of one run.\n$lastHexCell = $table.Cell($rowIndex, 2)\nSet-CellRuns $lastHexCell @(\"7A67\")\n\n$lastBinaryCell = $table.Cell($rowIndex, 3)\nSet-CellRuns $lastBinaryCell @(\"0111\", \" \", \"1010\", \" \", \"0110\", \" \", \"0111\")\n"}
